$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (some look numeric, e.g. "132.40"); force
# the cells to Text format first so Excel does not silently normalize them
# (stripping insignificant trailing zeros) when the .Value is assigned.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.080.50"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "1.815.56"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "337.71"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "0.9984"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "0.4149"
$ws.Range("E7").Value = "  +8.41%  "
$ws.Range("D8").Value = "0.3482"
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("D9").Value = "45.74"
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("D10").Value = "1.155"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").Value = "0.07480"
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").Value = "22.87"
$ws.Range("E12").Value = "  -3.11%  "
$ws.Range("D13").Value = "1.001"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "6.285"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D15").Value = "7.301"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").Value = "1.810.66"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "0.06672"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "82.34"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "0.9987"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").Value = "17.34"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").Value = "6.375"
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("D23").Value = "28.145.94"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("E24").Value = "  -1.91%  "
$ws.Range("D25").Value = "2.414"
$ws.Range("E25").Value = "  +2.10%  "
$ws.Range("D26").Value = "2.468"
$ws.Range("E26").Value = "  +2.50%  "
$ws.Range("D28").Value = "155.59"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").Value = "2.015.82"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").Value = "1.323"
$ws.Range("E30").Value = "  -7.03%  "
$ws.Range("D31").Value = "132.40"
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("D32").Value = "4.078"
$ws.Range("E32").Value = "  +1.43%  "
$ws.Range("D33").Value = "6.010"
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("D34").Value = "0.09033"
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("D35").Value = "12.38"
$ws.Range("E35").Value = "  -2.99%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.06327"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.02345"
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("D38").Value = "0.6667"
$ws.Range("E38").Value = "  -2.74%  "
$ws.Range("D39").Value = "5.243"
$ws.Range("E39").Value = "  -2.13%  "
$ws.Range("D40").Value = "0.2157"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("D41").Value = "1.516"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("E42").Value = "  -2.71%  "
$ws.Range("D43").Value = "8.146"
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("D44").Value = "14.24"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "0.9987"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("D46").Value = "0.6164"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("D47").Value = "3.878"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "128.55"
$ws.Range("E48").Value = "  -3.60%  "
$ws.Range("D49").Value = "2.054"
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("D50").Value = "1.181"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").Value = "0.07117"
$ws.Range("E51").Value = "  -5.29%  "
